# Applies a cyclic rotation of the species/observation data held in
# columns A,B,D,E,F,G,H,Q,R (and the optional "Ringhack" comment in AC)
# across rows 2-7, while leaving the per-row location/metadata columns
# (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY)
# untouched. This mirrors the upstream commit's row-by-row value swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row after the edit (taken from the canonical diff).
$rowsAfter = @{
    2 = @{ A = 89596652; B = 89392;  D = "NT"; E = 1202;   F = "Ullticka";           G = "Phellinidium ferrugineofuscum"; H = "(P.Karst.) Fiasson & Niemelä"; Q = 499487.0027453392; R = 7085828.909325912; AC = $null }
    3 = @{ A = 89596664; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";            G = "Alectoria sarmentosa";           H = "(Ach.) Ach.";                   Q = 499509.930857877;  R = 7085817.01519425;  AC = $null }
    4 = @{ A = 89596650; B = 56395;  D = "NT"; E = 100109; F = "Tretåig hackspett";  G = "Picoides tridactylus";           H = "(Linnaeus, 1758)";              Q = 499359.1099011709; R = 7085759.791431802; AC = "Ringhack" }
    5 = @{ A = 89596702; B = 89356;  D = "LC"; E = 5447;   F = "Vedticka";           G = "Fuscoporia viticola";            H = "(Schwein.) Murrill";            Q = 499326.9311505322; R = 7085816.168399233; AC = $null }
    6 = @{ A = 89596704; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";            G = "Alectoria sarmentosa";           H = "(Ach.) Ach.";                   Q = 498982.9346804961; R = 7085676.214588321; AC = $null }
    7 = @{ A = 89596600; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";            G = "Alectoria sarmentosa";           H = "(Ach.) Ach.";                   Q = 499800.9774889108; R = 7085994.014255985; AC = $null }
}

foreach ($r in 2..7) {
    $vals = $rowsAfter[$r]

    $ws.Range("A$r").Value = $vals.A
    $ws.Range("B$r").Value = $vals.B
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R

    if ($vals.AC) {
        $ws.Range("AC$r").Value = $vals.AC
    } else {
        $ws.Range("AC$r").ClearContents()
    }
}
